$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled data update).
# Column D holds price text that sometimes looks like a plain number
# (e.g. "246.37", "12.40"); a leading apostrophe forces Excel to keep it
# as literal text (preserving trailing zeros / exact formatting) instead
# of silently converting it to a floating point value.
$ws.Range("D2").Value = '35.459.50'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.903.32'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '''246.37'
$ws.Range("E5").Value = '  +3.26%  '
$ws.Range("D6").Value = '''0.646'
$ws.Range("E6").Value = '  +3.81%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '''41.89'
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("D9").Value = '''0.342'
$ws.Range("E9").Value = '  +3.90%  '
$ws.Range("D10").Value = '''0.0706'
$ws.Range("E10").Value = '  +1.37%  '
$ws.Range("D11").Value = '''0.0998'
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '2.179.81'
$ws.Range("E12").Value = '  +2.40%  '
$ws.Range("D13").Value = '''12.40'
$ws.Range("E13").Value = '  +8.82%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.912.24'
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.696'
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").Value = '''4.82'
$ws.Range("E16").Value = '  +2.46%  '
$ws.Range("D17").Value = '35.504.72'
$ws.Range("E17").Value = '  +1.27%  '
$ws.Range("D18").Value = '''71.97'
$ws.Range("E18").Value = '  +2.35%  '
$ws.Range("D19").Value = '0.0₃0828'
$ws.Range("E19").Value = '  +3.87%  '
$ws.Range("D20").Value = '''242.96'
$ws.Range("E20").Value = '  +0.81%  '
$ws.Range("D21").Value = '''12.61'
$ws.Range("E21").Value = '  +3.36%  '
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("E25").Value = '  +14.94%  '
$ws.Range("D26").Value = '''171.72'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '''8.52'
$ws.Range("E27").Value = '  +7.58%  '
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '''0.127'
$ws.Range("E29").Value = '  +1.94%  '
$ws.Range("D30").Value = '''0.964'
$ws.Range("E30").Value = '  +23.74%  '
$ws.Range("E31").Value = '  +1.88%  '
$ws.Range("E32").Value = '  +3.13%  '
$ws.Range("E33").Value = '  +0.17%  '
$ws.Range("E34").Value = '  +4.08%  '
$ws.Range("E35").Value = '  +8.49%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '''1.34'
$ws.Range("E37").Value = '  +2.69%  '
$ws.Range("D38").Value = '''1.11'
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("E39").Value = '  +16.99%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("D41").Value = '''91.60'
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''15.67'
$ws.Range("E42").Value = '  +5.29%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.345.98'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = '''49.04'
$ws.Range("E44").Value = '  +41.15%  '
$ws.Range("D45").Value = '''2.39'
$ws.Range("E45").Value = '  +2.12%  '
$ws.Range("D46").Value = '''12.73'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("D49").Value = '''6.59'
$ws.Range("E49").Value = '  +3.64%  '
$ws.Range("D50").Value = '2.091.24'
$ws.Range("E50").Value = '  +2.36%  '
$ws.Range("D51").Value = '''0.0693'
$ws.Range("E51").Value = '  +1.85%  '
